# Update modelxgb_df_results.xlsx with refreshed metrics (R^2 / RMSE / U)
# and the corresponding heat-map style shading for the RMSE (D) and U (E)
# columns, matching the freshly re-run notebook output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ConvertHexToBgr($hexColorValue) {
    $compR = [Convert]::ToInt32($hexColorValue.Substring(0, 2), 16)
    $compG = [Convert]::ToInt32($hexColorValue.Substring(2, 2), 16)
    $compB = [Convert]::ToInt32($hexColorValue.Substring(4, 2), 16)
    return ($compB * 65536) + ($compG * 256) + $compR
}

# rowNum => R^2, RMSE, U  (columns C, D, E)
$metricsByRow = @{
    2  = @(0.2438, 0.0401, 1.7724)
    3  = @(0.4039, 0.0354, 1.5715)
    4  = @(0.6251, 0.028,  1.2076)
    5  = @(0.8027, 0.0203, 0.8639)
    6  = @(0.9858, 0.0055, 0.24)
    7  = @(0.7663, 0.0222, 0.9809)
    8  = @(0.5944, 0.0295, 1.3077)
    9  = @(0.3822, 0.0367, 1.6166)
    10 = @(0.2225, 0.0415, 1.8331)
}

# rowNum => fill colors for RMSE (D) and U (E) cells
$fillsByRow = @{
    2  = @("F2FAEF", "F2FAEF")
    3  = @("DBF1D5", "DBF1D6")
    4  = @("A0D99B", "9BD696")
    5  = @("50B264", "48AE60")
    6  = @("00441B", "00441B")
    7  = @("65BD6F", "66BD6F")
    8  = @("AEDEA7", "AFDFA8")
    9  = @("E3F4DE", "E3F4DE")
    10 = @("F7FCF5", "F7FCF5")
}

foreach ($rowIndex in 2..10) {
    $rowMetrics = $metricsByRow[$rowIndex]
    $ws.Cells.Item($rowIndex, 3).Value = $rowMetrics[0]
    $ws.Cells.Item($rowIndex, 4).Value = $rowMetrics[1]
    $ws.Cells.Item($rowIndex, 5).Value = $rowMetrics[2]

    $rowFills = $fillsByRow[$rowIndex]
    $ws.Cells.Item($rowIndex, 4).Interior.Color = ConvertHexToBgr($rowFills[0])
    $ws.Cells.Item($rowIndex, 5).Interior.Color = ConvertHexToBgr($rowFills[1])
}
